# "Connected Office Test Data.xlsx" - final commit
# The "Test Results" sheet's header row (ID / Create Test Passed / Read Test
# Passed / Update Test Passed / Delete Test Passed) is removed and every
# row's pass/fail flags flip from FALSE to TRUE (all CRUD tests now pass).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")
$ws.Activate()

# New column-A id sequence once the header row is gone (rows shift up); the
# final row re-uses the last id since nothing shifts in behind it.
$ids = @("Z01","Z02","Z03","Z04","Z05","Z06","Z07","Z08","Z09", `
         "C01","C02","C03","C04","C05", `
         "D01","D02","D03","D04","D05","D06","D07","D08","D09","D09")

for ($i = 0; $i -lt $ids.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $ids[$i]
    $ws.Cells.Item($r, 2).Value = $true
    $ws.Cells.Item($r, 3).Value = $true
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = $true
}

# Scroll the viewport down and leave the final cell selected, matching
# where the user's cursor ended up after reviewing the last test result.
$ws.Range("E24").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
